$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Consumer Price Index, 1800-" and "Mountain Meadows Association" rows
# (previously the last two data rows) move up to become rows 3 and 4,
# pushing "Deseret News..." and "Mountain Meadows Monument Foundation" down
# to rows 5 and 6. Columns B/C/D are identical across these rows, so only
# column A (title) and column E (uri + hyperlink) actually need updating.

$ws.Range("A3").Value = "Consumer Price Index, 1800-"
$ws.Range("A4").Value = "Mountain Meadows Association"
$ws.Range("A5").Value = "Deseret News | 1875-08-04 | Brigham Young"
$ws.Range("A6").Value = "Mountain Meadows Monument Foundation"

$ws.Range("E3").Value = "https://www.minneapolisfed.org/about-us/monetary-policy/inflation-calculator/consumer-price-index-1800-"
$ws.Range("E4").Value = "http://www.mtn-meadows-assoc.com/"
$ws.Range("E5").Value = "https://newspapers.lib.utah.edu/details?id=2641490"
$ws.Range("E6").Value = "http://www.mmmf.org/"

# Rebuild the hyperlinks to match the new uri values per row (the engine
# re-adds them against the collection as a whole).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E2"), "https://secure.flickr.com/photos/jstephenconn/2807773224/")
$ws.Hyperlinks.Add($ws.Range("E3"), "https://www.minneapolisfed.org/about-us/monetary-policy/inflation-calculator/consumer-price-index-1800-")
$ws.Hyperlinks.Add($ws.Range("E4"), "http://www.mtn-meadows-assoc.com/")
$ws.Hyperlinks.Add($ws.Range("E5"), "https://newspapers.lib.utah.edu/details?id=2641490")
$ws.Hyperlinks.Add($ws.Range("E6"), "http://www.mmmf.org/")

# Re-adding hyperlinks stamps a freshly duplicated "Hyperlink" style onto the
# cell; reapply the named style so the cells keep using the workbook's
# original Hyperlink cell style instead of an extra near-duplicate one.
$ws.Range("E2").Style = "Hyperlink"
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("E5").Style = "Hyperlink"
$ws.Range("E6").Style = "Hyperlink"
